$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "TTT move cursor" was renamed / refined to "TTT move cursor into edge"
$ws.Range("C6").Value = "TTT move cursor into edge"

# Mark the rows that are now implemented/approved with the built-in "Good"
# style (green). Row 8 (TTT complete match) and the two footer rows
# (14/15) stay in their previous "Bad" (red) style.
$ws.Range("D2:D7").Style = "Good"
$ws.Range("D9:D13").Style = "Good"

# Narrow the Approved/Implemented columns now that they just hold a
# colored status flag rather than free text.
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 10.833333333333334

# Leave the cursor on F18 like the author did before saving.
$ws.Range("F18").Select() | Out-Null
